$wb = $excel.ActiveWorkbook

# Update version number on the isa_template sheet (metadata sheet)
$wsMeta = $wb.Worksheets.Item("isa_template")
$wsMeta.Range("B4").Value = "1.0.1"

# Add example values to the mandatory ENA annotation table on "New Table" sheet
$wsData = $wb.Worksheets.Item("New Table")

$wsData.Range("B2").Value = "Zea mays"
$wsData.Range("C2").Value = "NCBITaxon"
$wsData.Range("D2").Value = "http://purl.obolibrary.org/obo/NCBITaxon_4577"
$wsData.Range("E2").Value = "no"
$wsData.Range("H2").Value = "RdRp"
$wsData.Range("K2").Value = "RNA dependent RNA polymerase"
$wsData.Range("L2").Value = "GO"
$wsData.Range("M2").Value = "http://purl.obolibrary.org/obo/GO_0003968"
$wsData.Range("N2").Value = "1"
$wsData.Range("Q2").Value = "1"
$wsData.Range("T2").Value = "3600"
$wsData.Range("W2").Value = "no"
$wsData.Range("Z2").Value = "no"
